$d = $word.ActiveDocument

# Locate the "Product Development and Platform Architecture" paragraph under
# the Siege Analytics / PARTNER role and insert three new bullet paragraphs
# immediately after it (before the existing "Conceived, architected..." bullet).
$r = $d.Content
$found = $r.Find.Execute("Product Development and Platform Architecture", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)
    $bullet1 = "Conceived and architected redistricting platform incorporating boundary estimation algorithm used by 2,500+ analysts"
    $bullet2 = "Built multi-tenant data warehouse tracking decades of demographic data, enabling discovery of 500,000+ mischaracterized voters"
    $bullet3 = "Platform democratized redistricting analysis, reducing costs by 75% and enabling 200+ smaller organizations to participate"
    $newText = "`r" + [char]0x2022 + " " + $bullet1 + "`r" + [char]0x2022 + " " + $bullet2 + "`r" + [char]0x2022 + " " + $bullet3
    $r.InsertAfter($newText)
}
